# Base Simulador - maio v2
#
# The underlying KPI simulation was re-run, which (a) regenerated most of
# the numeric KPI columns (small floating point drift throughout, larger
# swings in KPI4/KPI6/KPI8 for several scenarios), (b) added previously
# missing KPI5 ("F" column) samples for many rows, and (c) produced one
# fewer scenario row overall (67 data rows instead of 68, old row 68 is
# gone and the former row-54 record was dropped/absorbed during the
# reshuffle). The header row (row 1) and column layout are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The re-simulation produced one less scenario than before, so the last
# data row (68) no longer exists - drop it and let the rows below it
# (none, in this case) shift up.
$ws.Rows.Item(68).Delete()

# Full refreshed KPI table for rows 2-67 (A:J), values taken from the
# updated simulation output.
$data = @(
    @([double]"0", [double]"247105.769379", [double]"352421.8203048419", [double]"0.999999995", [double]"33686531.32457296", [double]"5", [double]"93328618.23009101", [double]"5", [double]"128432757.136888", "S10"),
    @([double]"1", [double]"217838.694291", [double]"335995.3795279213", [double]"0.999999995", [double]"-4332350.325584246", [double]"26", [double]"103079563.372284", [double]"5", [double]"88394348.4730044", "S10"),
    @([double]"2", [double]"266863.3809499999", [double]"295364.757463744", [double]"-7E-09", [double]"6872304.585514694", [double]"19", [double]"60776319.5446833", [double]"5", [double]"53121362.2945007", "S10"),
    @([double]"3", [double]"84468.887581", [double]"228429.1621248172", [double]"5.999999985", [double]"13833365.16757776", [double]"12", [double]"81070536.6320347", [double]"5", [double]"65436270.2644522", "S10"),
    @([double]"4", [double]"49385.5195", [double]"193850.697720596", [double]"1.999999997", [double]"-395364.8421912249", [double]"31", [double]"36950421.4232691", [double]"5", [double]"24156638.5693997", "S10"),
    @([double]"5", [double]"279246.515913", $null, [double]"4.999999985", [double]"28271438.20374005", [double]"2", $null, [double]"5", [double]"132559336.26135", "S10"),
    @([double]"6", [double]"118411.95", [double]"268796.303622337", [double]"2.999999996", [double]"-1277732.629437586", [double]"24", [double]"65867510.0509068", [double]"5", [double]"49886784.6835633", "S10"),
    @([double]"7", [double]"226907.941727", [double]"313322.1007418785", [double]"0.999999992", [double]"4720467.70784786", [double]"23", [double]"89465155.19588719", [double]"5", [double]"70133146.86166421", "S10"),
    @([double]"8", [double]"27423.9139", [double]"165757.7205671078", [double]"1.999999997", [double]"9447875.952797845", [double]"28", [double]"66879630.3622337", [double]"5", [double]"48056588.2187911", "S10"),
    @([double]"9", [double]"224764.57836", [double]"283070.9132577037", [double]"-7E-09", [double]"12137737.88375509", [double]"16", [double]"66004327.5546375", [double]"5", [double]"62031125.3162258", "S10"),
    @([double]"10", [double]"311633.7437350002", [double]"627191.375118786", [double]"-6E-09", [double]"20709427.39561249", [double]"13", [double]"120520866.772065", [double]"5", [double]"111336706.643765", "S10"),
    @([double]"11", [double]"103907.158792", [double]"304635.8107719819", [double]"0.999999986", [double]"9546561.498471471", [double]"30", [double]"95636875.7250993", [double]"5", [double]"78416605.398717", "S10"),
    @([double]"12", [double]"136114.660307", [double]"283396.3395450573", [double]"2.99999999", [double]"163899.5297748353", [double]"25", [double]"88162403.00945599", [double]"5", [double]"58912880.5055408", "S10"),
    @([double]"13", [double]"0", [double]"188157.4125356359", [double]"3.999999995", [double]"3416998.1415625", [double]"29", [double]"44799383.9370562", [double]"5", [double]"8130155.51005326", "S10"),
    @([double]"14", [double]"112651.161581", [double]"268796.303622337", [double]"4.99999999", [double]"-688448.9842057356", [double]"21", [double]"57183396.5740578", [double]"5", [double]"44025324.5983992", "S10"),
    @([double]"15", [double]"250751.614062", [double]"594753.0371976423", [double]"0.999999993", [double]"8032350.7650893", [double]"22", [double]"140443579.463008", [double]"5", [double]"123566046.138903", "S10"),
    @([double]"16", [double]"166713.6957030001", [double]"459745.1408452092", [double]"1.999999995", [double]"10540497.57212998", [double]"17", [double]"113929898.858352", [double]"5", [double]"116424007.92435", "S10"),
    @([double]"17", [double]"176563.898581", [double]"245810.8362400445", [double]"1.999999993", [double]"9630383.649145005", [double]"9", [double]"57316305.532588", [double]"5", [double]"71559483.5310411", "S10"),
    @([double]"18", [double]"118179.140522", [double]"141655.7045225727", [double]"0.999999996", [double]"12709479.21787969", [double]"8", [double]"32909016.7893645", [double]"5", [double]"39169642.8058211", "S10"),
    @([double]"19", [double]"98478.38504199999", [double]"197440.5254561626", [double]"1.999999992", [double]"13483984.57258214", [double]"10", [double]"43264457.1595386", [double]"5", [double]"50942225.2829144", "S10"),
    @([double]"20", [double]"82825.84376", [double]"268796.303622337", [double]"0.9999999939999999", [double]"8735829.482966224", [double]"27", [double]"69444871.20746569", [double]"5", [double]"66609639.7023026", "S10"),
    @([double]"21", [double]"204415.709125", [double]"268796.303622337", [double]"2.999999989", [double]"6345431.665514652", [double]"11", [double]"68765064.9344134", [double]"5", [double]"61783181.6933618", "S10"),
    @([double]"22", [double]"194163.913047", [double]"436224.5932164576", [double]"2.999999993", [double]"4440189.221211218", [double]"15", [double]"103009033.986665", [double]"5", [double]"95985094.22961649", "S10"),
    @([double]"23", [double]"565179.4984789999", [double]"564472.2376069077", [double]"14.999999984", [double]"79762159.1690679", [double]"1", [double]"118207994.076012", [double]"5", [double]"160152293.874123", "S10"),
    @([double]"24", [double]"210427.546306", [double]"564472.2376069077", [double]"7.999999991", [double]"12697706.55173514", [double]"4", [double]"125055714.234396", [double]"5", [double]"125474356.64008", "S10"),
    @([double]"25", [double]"201118.835714", [double]"447993.8393705615", [double]"3.999999989", [double]"-5273485.49575359", [double]"20", [double]"93543166.2384655", [double]"5", [double]"75236861.78298929", "S10"),
    @([double]"26", [double]"309371.4254999999", [double]"627191.375118786", [double]"2.999999987", [double]"67157403.89008105", [double]"3", [double]"128532347.967239", [double]"5", [double]"168993135.734353", "S10"),
    @([double]"27", [double]"149959.284461", [double]"179197.5357482248", [double]"0.999999995", [double]"9585279.296527943", [double]"18", [double]"35664507.8665033", [double]"5", [double]"29011722.1974944", "S10"),
    @([double]"28", [double]"3202.44", [double]"179491.3867783297", [double]"-1E-09", [double]"749228.672174304", [double]"32", [double]"23540432.1257414", [double]"5", [double]"3769167.73774433", "S10"),
    @([double]"29", [double]"113004.929588", $null, [double]"10.999999979", [double]"11784528.38437979", [double]"7", $null, [double]"5", [double]"54015049.9983959", "S10"),
    @([double]"30", [double]"61390.669431", [double]"161277.7821734023", [double]"3.999999995", [double]"49143446.77005175", [double]"6", [double]"40319445.5433506", [double]"5", [double]"62035701.8863961", "S10"),
    @([double]"31", [double]"206845.6237900001", [double]"517761.8563709559", [double]"2.999999992", [double]"6069764.464888243", [double]"14", [double]"128307078.614489", [double]"5", [double]"131135133.148493", "S10"),
    @([double]"32", [double]"437041.5331989997", $null, [double]"4.999999985", [double]"14129465.69103084", [double]"22", $null, [double]"5", [double]"228903868.273275", "G10"),
    @([double]"33", [double]"1427756.390602", [double]"1989956.15493068", [double]"-9E-09", [double]"-7245394.153791296", [double]"17", [double]"462914368.036666", [double]"5", [double]"419200040.408371", "G10"),
    @([double]"34", [double]"160878.3722619999", [double]"401593.7734454345", [double]"4.999999993", [double]"15674217.66116987", [double]"30", [double]"162971737.218567", [double]"5", [double]"140937442.588491", "G10"),
    @([double]"35", [double]"547208.544608", [double]"716704.011697341", [double]"1.99999999", [double]"19674320.38476176", [double]"13", [double]"165949842.129764", [double]"5", [double]"165646739.058223", "G10"),
    @([double]"36", [double]"659069.2998729999", [double]"825167.6134377619", [double]"7.999999976", [double]"27574170.14620422", [double]"5", [double]"195397061.094805", [double]"5", [double]"200151254.904617", "G10"),
    @([double]"37", [double]"1166441.761158001", $null, [double]"-2E-08", [double]"9073580.796574241", [double]"16", $null, [double]"5", [double]"348565330.671247", "G10"),
    @([double]"38", [double]"324307.9065659997", [double]"582391.9911817291", [double]"0.9999999939999999", [double]"2744269.59513884", [double]"27", [double]"264508335.976307", [double]"5", [double]"251447120.831897", "G10"),
    @([double]"39", [double]"261352.09457", [double]"648949.418757872", [double]"-4.000000007", [double]"22765456.63839758", [double]"32", [double]"180773023.234012", [double]"5", [double]"130811153.108373", "G10"),
    @([double]"40", [double]"485061.0454419998", $null, [double]"1.999999985", [double]"21084562.53407755", [double]"6", $null, [double]"5", [double]"152146785.390455", "G10"),
    @([double]"41", [double]"456089.5878179997", [double]"1143139.510627199", [double]"-2.000000014", [double]"8365578.784704199", [double]"31", [double]"269938285.450157", [double]"5", [double]"240470886.97798", "G10"),
    @([double]"42", [double]"1045518.143118", [double]"1870576.611309708", [double]"8.999999957", [double]"67699177.56819198", [double]"8", [double]"453930494.313569", [double]"5", [double]"398308822.759385", "G10"),
    @([double]"43", [double]"595891.9738200001", [double]"725750.019780311", [double]"-1.000000039", [double]"31097026.64363621", [double]"7", [double]"338347044.616148", [double]"5", [double]"346315867.113864", "G10"),
    @([double]"44", [double]"901451.858009999", [double]"2943319.52466459", [double]"4.999999985", [double]"37792529.97965543", [double]"24", [double]"871844868.1560791", [double]"5", [double]"356571751.166594", "G10"),
    @([double]"45", [double]"332567.650706", [double]"447993.8393705615", [double]"-1.000000008", [double]"1489812.729699546", [double]"28", [double]"392707751.013488", [double]"5", [double]"351708377.004988", "G10"),
    @([double]"46", [double]"389937.957858", [double]"685956.607590682", [double]"0.999999991", [double]"4327483.929702809", [double]"25", [double]"127932972.017863", [double]"5", [double]"139025695.080742", "G10"),
    @([double]"47", [double]"330437.077012", [double]"782269.404911224", [double]"0.9999999899999999", [double]"9157319.962886602", [double]"29", [double]"188896780.837437", [double]"5", [double]"173186776.487825", "G10"),
    @([double]"48", [double]"704143.3872839999", $null, [double]"11.999999967", [double]"99004889.16466612", [double]"10", $null, [double]"5", [double]"375049738.363188", "G10"),
    @([double]"49", [double]"738857.865405001", [double]"1495985.827810116", [double]"-1.4E-08", [double]"-4070848.855628334", [double]"26", [double]"353258554.43773", [double]"5", [double]"294729975.30405", "G10"),
    @([double]"50", [double]"399500.4786539997", [double]"1479373.320258577", [double]"0.999999985", [double]"2406525.90925815", [double]"33", [double]"349335762.370402", [double]"5", [double]"285022400.390091", "G10"),
    @([double]"51", [double]"1664608.373266002", [double]"2343958.091319741", [double]"5.999999982", [double]"36880788.4355883", [double]"2", [double]"490132517.48139", [double]"5", [double]"475500342.421061", "G10"),
    @([double]"52", [double]"1192790.478853999", $null, [double]"-5.000000029", [double]"24041690.26878395", [double]"21", $null, [double]"5", [double]"499237175.937533", "G10"),
    @([double]"53", [double]"295237.1833050001", [double]"766539.41086116", [double]"0.999999991", [double]"-2412282.298719347", [double]"34", [double]"181008866.249681", [double]"5", [double]"150776833.048063", "G10"),
    @([double]"54", [double]"559147.208685", [double]"671990.759055843", [double]"4.999999989", [double]"46342348.84790082", [double]"9", [double]"280524877.955404", [double]"5", [double]"238815037.368466", "G10"),
    @([double]"55", [double]"948751.756936", [double]"1075185.214489349", [double]"2.999999982", [double]"19574394.36861564", [double]"4", [double]"443866203.113601", [double]"5", [double]"402289284.770168", "G10"),
    @([double]"56", [double]"1030383.516133", [double]"1075185.214489349", [double]"2.999999985", [double]"-18291582.73189784", [double]"11", [double]"435872490.965229", [double]"5", [double]"372284378.690909", "G10"),
    @([double]"57", [double]"867290.346103", [double]"1075185.214489349", [double]"-1.1E-08", [double]"-16038339.70004725", [double]"20", [double]"534243142.793327", [double]"5", [double]"464730008.847567", "G10"),
    @([double]"58", [double]"830777.7814730001", [double]"1075185.214489349", [double]"-1.2E-08", [double]"-7755943.792864057", [double]"15", [double]"347512645.701532", [double]"5", [double]"403819992.292377", "G10"),
    @([double]"59", [double]"847785.775594", [double]"1075185.214489349", [double]"-2.000000013", [double]"14885126.86672064", [double]"12", [double]"405377230.617039", [double]"5", [double]"390448045.247282", "G10"),
    @([double]"60", [double]"463968.8240989997", [double]"761589.526929955", [double]"4.999999982", [double]"16970925.28966692", [double]"19", [double]"310420875.883594", [double]"5", [double]"313137532.893437", "G10"),
    @([double]"61", [double]"431734.2034979999", [double]"456771.5319966288", [double]"1.999999992", [double]"5398299.655918254", [double]"18", [double]"143247730.662296", [double]"5", [double]"127360244.329447", "G10"),
    @([double]"62", [double]"3632141.897297002", $null, [double]"4.999999965", [double]"171051883.5710644", [double]"1", $null, [double]"5", [double]"1059393350.79388", "G10"),
    @([double]"63", [double]"1404228.462323", $null, [double]"0.999999976", [double]"16745936.93818602", [double]"14", $null, [double]"5", [double]"513017299.883275", "G10"),
    @([double]"64", [double]"708795.817177", [double]"940787.0626781799", [double]"-1.1E-08", [double]"-4735994.28343372", [double]"23", [double]"642220644.6680959", [double]"5", [double]"602016352.158518", "G10"),
    @([double]"65", [double]"1117577.921944003", [double]"1363203.589777557", [double]"2.99999998", [double]"20982985.68700553", [double]"3", [double]"321903714.794352", [double]"5", [double]"303764449.160994", "G10")
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
